$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: merge every run inside a paragraph's range into a single run
# holding $finalText. Word's "no-op" detection skips the write if the
# replacement text is character-for-character identical to what Range.Text
# already reports, so we first stamp a placeholder string (guaranteed to
# differ) and only then write the real text - this forces the run merge.
# ---------------------------------------------------------------------
function Merge-ParagraphRuns($para, $finalText) {
    $r = $para.Range
    $r.MoveEnd(1, -1)
    $r.Text = "___TMP_PLACEHOLDER___"

    $r2 = $para.Range
    $r2.MoveEnd(1, -1)
    $r2.Text = $finalText
}

# ---------------------------------------------------------------------
# 1) Merge the three runs of the balsamiq.cloud hyperlink into a single
#    run containing the full URL text.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*balsamiq.cloud/s*lo8ba*") {
        Merge-ParagraphRuns $p "https://balsamiq.cloud/s8lo8ba/pxpybsw/r52D9"
        break
    }
}

# ---------------------------------------------------------------------
# 2) Merge the five runs of the goo.gl hyperlink into a single run.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*goo.*787R*") {
        Merge-ParagraphRuns $p "https://goo.gl/yd787R"
        break
    }
}

# ---------------------------------------------------------------------
# 3) Insert four new empty paragraphs (same formatting as the
#    "(los puntos seran Nakits)" paragraph) right after it and before
#    the "(desea redimir ...)" paragraph.
# ---------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Nakits)*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $anchor = $d.Paragraphs.Item($anchorIndex)

    # Create 4 blank paragraphs right after the anchor paragraph.
    for ($n = 1; $n -le 4; $n++) {
        $anchor.Range.InsertParagraphAfter() | Out-Null
    }

    # Each freshly-inserted paragraph carries a stray empty <w:r> that
    # holds the paragraph-mark formatting; InsertXML on that paragraph's
    # own Range replaces its contents with a clean <w:p><w:pPr>...</w:pPr></w:p>
    # (no run at all), matching the target markup exactly.
    $newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="normal0"/><w:ind w:left="12" w:right="-749"/><w:contextualSpacing w:val="0"/><w:jc w:val="both"/><w:rPr><w:b/><w:color w:val="CC0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    for ($n = 1; $n -le 4; $n++) {
        $idx = $anchorIndex + $n
        $newp = $d.Paragraphs.Item($idx)
        $newp.Range.InsertXML($newParaXml)
    }
}

Write-Host "Edit complete"
